# Add the new "QBEA Q-Values" worksheet as the last tab (becomes the active sheet,
# which also clears tabSelected on the other sheets and sets workbookView activeTab).
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "QBEA Q-Values"

# Prime the shared-string table so new strings are interned in the same order
# as the target workbook (two-space placeholder must become the first new string).
$ws.Range("A2").Value = "  "

# Row 1
$ws.Range("A1").Value = "Running experiment QLearning bias: 100 repeated 1000 averaged over 1"

# Row 2
$ws.Range("A2").Value = "  "
$ws.Range("B2").Value = -10
$ws.Range("C2").Value = -9
$ws.Range("D2").Value = -8
$ws.Range("E2").Value = -7
$ws.Range("F2").Value = -6
$ws.Range("G2").Value = -5
$ws.Range("H2").Value = -4
$ws.Range("I2").Value = -3
$ws.Range("J2").Value = -2
$ws.Range("K2").Value = -1
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 7
$ws.Range("T2").Value = 8
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 10

# Row 3
$ws.Range("A3").Value = -5
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 253.94
$ws.Range("D3").Value = 227.79
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 223.67
$ws.Range("G3").Value = 57.99
$ws.Range("H3").Value = 143.23
$ws.Range("I3").Value = 57.03
$ws.Range("J3").Value = 182.86
$ws.Range("K3").Value = 237.64
$ws.Range("L3").Value = 257.21
$ws.Range("M3").Value = 175.28
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 125.28
$ws.Range("P3").Value = 258.45
$ws.Range("Q3").Value = 172.87
$ws.Range("R3").Value = 249.84
$ws.Range("S3").Value = 214.91
$ws.Range("T3").Value = 271.81
$ws.Range("U3").Value = 335.57
$ws.Range("V3").Value = 212.57

# Row 4
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 122.2
$ws.Range("C4").Value = 172.7
$ws.Range("D4").Value = 257.38
$ws.Range("E4").Value = 257.34
$ws.Range("F4").Value = 85.72
$ws.Range("G4").Value = 154.35
$ws.Range("H4").Value = 163.14
$ws.Range("I4").Value = 128.42
$ws.Range("J4").Value = 236.73
$ws.Range("K4").Value = 273.25
$ws.Range("L4").Value = 112.09
$ws.Range("M4").Value = 272.23
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 155.28
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 258.64
$ws.Range("R4").Value = 344.9
$ws.Range("S4").Value = 106.96
$ws.Range("T4").Value = 247.98
$ws.Range("U4").Value = 33.21
$ws.Range("V4").Value = 141.8

# Row 5
$ws.Range("A5").Value = "Averaged reward:"
$ws.Range("B5").Value = 102200
$ws.Range("C5").Value = "QLearning"

# Row 6
$ws.Range("A6").Value = "Time: 0.043 sec."

# Row 7
$ws.Range("A7").Value = "Running experiment QBEA bias: 100 repeated 1000 averaged over 1"

# Row 8
$ws.Range("A8").Value = "  "
$ws.Range("B8").Value = -10
$ws.Range("C8").Value = -9
$ws.Range("D8").Value = -8
$ws.Range("E8").Value = -7
$ws.Range("F8").Value = -6
$ws.Range("G8").Value = -5
$ws.Range("H8").Value = -4
$ws.Range("I8").Value = -3
$ws.Range("J8").Value = -2
$ws.Range("K8").Value = -1
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 3
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 5
$ws.Range("R8").Value = 6
$ws.Range("S8").Value = 7
$ws.Range("T8").Value = 8
$ws.Range("U8").Value = 9
$ws.Range("V8").Value = 10

# Row 9
$ws.Range("A9").Value = -5
$ws.Range("B9").Value = 145.11
$ws.Range("C9").Value = 147.14
$ws.Range("D9").Value = 149.11
$ws.Range("E9").Value = 151.11
$ws.Range("F9").Value = 153.1
$ws.Range("G9").Value = 156.41
$ws.Range("H9").Value = 153.09
$ws.Range("I9").Value = 151.27
$ws.Range("J9").Value = 149.15
$ws.Range("K9").Value = 147.11
$ws.Range("L9").Value = 145.11
$ws.Range("M9").Value = 143.11
$ws.Range("N9").Value = 141.11
$ws.Range("O9").Value = 139.11
$ws.Range("P9").Value = 137.11
$ws.Range("Q9").Value = 135.11
$ws.Range("R9").Value = 133.11
$ws.Range("S9").Value = 131.11
$ws.Range("T9").Value = 129.11
$ws.Range("U9").Value = 127.11
$ws.Range("V9").Value = 125.11

# Row 10
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = 7.29
$ws.Range("C10").Value = 9.29
$ws.Range("D10").Value = 11.29
$ws.Range("E10").Value = 13.29
$ws.Range("F10").Value = 15.29
$ws.Range("G10").Value = 17.29
$ws.Range("H10").Value = 19.29
$ws.Range("I10").Value = 21.29
$ws.Range("J10").Value = 23.29
$ws.Range("K10").Value = 25.29
$ws.Range("L10").Value = 27.29
$ws.Range("M10").Value = 29.29
$ws.Range("N10").Value = -58.79
$ws.Range("O10").Value = 33.29
$ws.Range("P10").Value = 35.28
$ws.Range("Q10").Value = 38.54
$ws.Range("R10").Value = 35.59
$ws.Range("S10").Value = 33.3
$ws.Range("T10").Value = 31.29
$ws.Range("U10").Value = 29.29
$ws.Range("V10").Value = 27.29

# Row 11
$ws.Range("A11").Value = "Averaged reward:"
$ws.Range("B11").Value = 25436
$ws.Range("C11").Value = "QBEA"

# Row 12
$ws.Range("A12").Value = "Time: 0.419 sec."

# Row 16
$ws.Range("A16").Value = "  "
$ws.Range("B16").Value = -10
$ws.Range("C16").Value = -9
$ws.Range("D16").Value = -8
$ws.Range("E16").Value = -7
$ws.Range("F16").Value = -6
$ws.Range("G16").Value = -5
$ws.Range("H16").Value = -4
$ws.Range("I16").Value = -3
$ws.Range("J16").Value = -2
$ws.Range("K16").Value = -1
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 1
$ws.Range("N16").Value = 2
$ws.Range("O16").Value = 3
$ws.Range("P16").Value = 4
$ws.Range("Q16").Value = 5
$ws.Range("R16").Value = 6
$ws.Range("S16").Value = 7
$ws.Range("T16").Value = 8
$ws.Range("U16").Value = 9
$ws.Range("V16").Value = 10

# Row 17
$ws.Range("A17").Value = -5
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 253.94
$ws.Range("D17").Value = 227.79
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 223.67
$ws.Range("G17").Value = 57.99
$ws.Range("H17").Value = 143.23
$ws.Range("I17").Value = 57.03
$ws.Range("J17").Value = 182.86
$ws.Range("K17").Value = 237.64
$ws.Range("L17").Value = 257.21
$ws.Range("M17").Value = 175.28
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 125.28
$ws.Range("P17").Value = 258.45
$ws.Range("Q17").Value = 172.87
$ws.Range("R17").Value = 249.84
$ws.Range("S17").Value = 214.91
$ws.Range("T17").Value = 271.81
$ws.Range("U17").Value = 335.57
$ws.Range("V17").Value = 212.57

# Row 18
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = 122.2
$ws.Range("C18").Value = 172.7
$ws.Range("D18").Value = 257.38
$ws.Range("E18").Value = 257.34
$ws.Range("F18").Value = 85.72
$ws.Range("G18").Value = 154.35
$ws.Range("H18").Value = 163.14
$ws.Range("I18").Value = 128.42
$ws.Range("J18").Value = 236.73
$ws.Range("K18").Value = 273.25
$ws.Range("L18").Value = 112.09
$ws.Range("M18").Value = 272.23
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 155.28
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 258.64
$ws.Range("R18").Value = 344.9
$ws.Range("S18").Value = 106.96
$ws.Range("T18").Value = 247.98
$ws.Range("U18").Value = 33.21
$ws.Range("V18").Value = 141.8

# Row 19
$ws.Range("A19").Value = "Averaged reward:"
$ws.Range("B19").Value = 102200
$ws.Range("C19").Value = "QLearning"

# Row 20
$ws.Range("A20").Value = "Time: 0.092 sec."

# Column widths (best effort through the character-width COM property; the engine
# snaps to a coarser grid than the bestFit pixel-perfect values in the source file).
$ws.Columns.Item(1).ColumnWidth = 55.5
$ws.Columns.Item(2).ColumnWidth = 6.16666666666667
$ws.Columns.Item(3).ColumnWidth = 8.33333333333333

# Selection + active sheet view state matching the source sheetView
$ws.Range("Q18").Select()
